$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Item(1).Delete()
$wb.Save()
